$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 71.55
$ws.Range("I9").Value = 67.8
$ws.Range("J9").Value = 82.8
$ws.Range("K9").Value = 67.8
$ws.Range("L9").Value = 82.8
$ws.Range("M9").Value = 101.2
$ws.Range("N9").Value = -420.8
$ws.Range("H33").Value = 716.381
$ws.Range("I33").Value = 244.4
$ws.Range("J33").Value = 1145.4546
$ws.Range("K33").Value = 244.4
$ws.Range("L33").Value = 1145.4546
$ws.Range("M33").Value = -15.40000000000001
$ws.Range("N33").Value = -1603.4546
$ws.Range("H41").Value = 1150.0
$ws.Range("I41").Value = 0.0
$ws.Range("J41").Value = 1150.0
$ws.Range("K41").Value = 0.0
$ws.Range("L41").Value = 1150.0
$ws.Range("M41").ClearContents()
$ws.Range("N41").Value = -2030.0
$ws.Range("H43").Value = 5132.077
$ws.Range("J43").Value = 8502.429
$ws.Range("L43").Value = 8502.429
$ws.Range("N43").Value = -8640.429
$ws.Range("H53").Value = 599.8947
$ws.Range("I53").Value = 467.18182
$ws.Range("K53").Value = 467.18182
$ws.Range("M53").Value = 169.81818
$ws.Range("H55").Value = 188.5
$ws.Range("I55").Value = 188.5
$ws.Range("K55").Value = 188.5
$ws.Range("M55").Value = 25.5
$ws.Range("H82").Value = 703.5
$ws.Range("I82").Value = 703.5
$ws.Range("K82").Value = 2110.5
$ws.Range("M82").Value = -1704.5
$ws.Range("H85").Value = 703.5
$ws.Range("I85").Value = 703.5
$ws.Range("K85").Value = 2110.5
$ws.Range("M85").Value = -706.5
$ws.Range("H93").Value = 0.0
$ws.Range("J93").Value = 0.0
$ws.Range("L93").Value = 0.0
$ws.Range("N93").ClearContents()
$ws.Range("H96").Value = 556283.25
$ws.Range("I96").Value = 833772.8
$ws.Range("K96").Value = 2501318.4
$ws.Range("M96").Value = -2499945.4
$ws.Range("H113").Value = 5600.0
$ws.Range("I113").Value = 4000.0
$ws.Range("J113").Value = 6000.0
$ws.Range("K113").Value = 4000.0
$ws.Range("L113").Value = 6000.0
$ws.Range("M113").Value = -746.0
$ws.Range("N113").Value = -12508.0
$ws.Range("H125").Value = 55196.25
$ws.Range("I125").Value = 0.0
$ws.Range("J125").Value = 55196.25
$ws.Range("K125").Value = 0.0
$ws.Range("L125").Value = 496766.25
$ws.Range("M125").ClearContents()
$ws.Range("N125").Value = -501686.25
$ws.Range("H132").Value = 2745.1155
$ws.Range("I132").Value = 2112.375
$ws.Range("K132").Value = 6337.125
$ws.Range("M132").Value = -3807.125
$ws.Range("H135").Value = 1631.9
$ws.Range("I135").Value = 1953.75
$ws.Range("J135").Value = 344.5
$ws.Range("K135").Value = 17583.75
$ws.Range("L135").Value = 3100.5
$ws.Range("M135").Value = -15048.75
$ws.Range("N135").Value = -8170.5
$ws.Range("H138").Value = 5085.87
$ws.Range("I138").Value = 1166.3334
$ws.Range("J138").Value = 5620.352
$ws.Range("K138").Value = 3499.0002
$ws.Range("L138").Value = 16861.056
$ws.Range("M138").Value = 1640.9998
$ws.Range("N138").Value = -27141.056

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2275.0303
$ws.Range("I2").Value = 2213.0
$ws.Range("J2").Value = 2359.2144
$ws.Range("K2").Value = 2213.0
$ws.Range("L2").Value = 2359.2144
$ws.Range("M2").Value = -2100.0
$ws.Range("N2").Value = -2585.2144
$ws.Range("H17").Value = 0.0
$ws.Range("J17").Value = 0.0
$ws.Range("L17").Value = 0.0
$ws.Range("N17").ClearContents()
$ws.Range("H35").Value = 2684.25
$ws.Range("I35").Value = 2684.25
$ws.Range("K35").Value = 2684.25
$ws.Range("M35").Value = -2278.25
$ws.Range("H50").Value = 7861.9
$ws.Range("I50").Value = 322.66666
$ws.Range("J50").Value = 11093.0
$ws.Range("K50").Value = 322.66666
$ws.Range("L50").Value = 11093.0
$ws.Range("M50").Value = 391.33334
$ws.Range("N50").Value = -12521.0
$ws.Range("H61").Value = 6081.316
$ws.Range("I61").Value = 5619.4707
$ws.Range("K61").Value = 5619.4707
$ws.Range("M61").Value = -5407.4707
$ws.Range("H80").Value = 50000.0
$ws.Range("J80").Value = 50000.0
$ws.Range("L80").Value = 50000.0
$ws.Range("N80").Value = -51996.0
$ws.Range("H83").Value = 50000.0
$ws.Range("J83").Value = 50000.0
$ws.Range("L83").Value = 150000.0
$ws.Range("N83").Value = -159984.0
$ws.Range("H116").Value = 2275.0303
$ws.Range("I116").Value = 2213.0
$ws.Range("J116").Value = 2359.2144
$ws.Range("K116").Value = 2213.0
$ws.Range("L116").Value = 2359.2144
$ws.Range("M116").Value = 81.0
$ws.Range("N116").Value = -6947.2144
$ws.Range("H132").Value = 2129.7576
$ws.Range("I132").Value = 2002.2069
$ws.Range("J132").Value = 3054.5
$ws.Range("K132").Value = 6006.620699999999
$ws.Range("L132").Value = 9163.5
$ws.Range("M132").Value = -3476.620699999999
$ws.Range("N132").Value = -14223.5
$ws.Range("H135").Value = 49999.5
$ws.Range("J135").Value = 49999.5
$ws.Range("L135").Value = 49999.5
$ws.Range("N135").Value = -60139.5
$ws.Range("H136").Value = 6081.316
$ws.Range("I136").Value = 5619.4707
$ws.Range("K136").Value = 16858.4121
$ws.Range("M136").Value = -14308.4121

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2275.0303
$ws.Range("I3").Value = 2213.0
$ws.Range("J3").Value = 2359.2144
$ws.Range("K3").Value = 2213.0
$ws.Range("L3").Value = 2359.2144
$ws.Range("M3").Value = -2099.0
$ws.Range("N3").Value = -2587.2144
$ws.Range("H7").Value = 764.0
$ws.Range("I7").Value = 57.2
$ws.Range("J7").Value = 1647.5
$ws.Range("K7").Value = 57.2
$ws.Range("L7").Value = 1647.5
$ws.Range("M7").Value = 55.8
$ws.Range("N7").Value = -1873.5
$ws.Range("H8").Value = 1690.2
$ws.Range("I8").Value = 1690.2
$ws.Range("K8").Value = 1690.2
$ws.Range("M8").Value = -1550.2
$ws.Range("H11").Value = 2090.5
$ws.Range("I11").Value = 636.0
$ws.Range("K11").Value = 636.0
$ws.Range("M11").Value = -496.0
$ws.Range("H12").Value = 106.5
$ws.Range("I12").Value = 78.0
$ws.Range("J12").Value = 249.0
$ws.Range("K12").Value = 78.0
$ws.Range("L12").Value = 249.0
$ws.Range("M12").Value = 90.0
$ws.Range("N12").Value = -585.0
$ws.Range("H14").Value = 2583.3333
$ws.Range("I14").Value = 1000.0
$ws.Range("K14").Value = 1000.0
$ws.Range("M14").Value = -828.0
$ws.Range("H20").Value = 1850.4286
$ws.Range("I20").Value = 2034.5555
$ws.Range("J20").Value = 1519.0
$ws.Range("K20").Value = 2034.5555
$ws.Range("L20").Value = 1519.0
$ws.Range("M20").Value = -1787.5555
$ws.Range("N20").Value = -2013.0
$ws.Range("H134").Value = 2500.0645
$ws.Range("I134").Value = 2042.3846
$ws.Range("K134").Value = 6127.1538
$ws.Range("M134").Value = -3592.1538
$ws.Range("H135").Value = 60780.0
$ws.Range("J135").Value = 60780.0
$ws.Range("L135").Value = 60780.0
$ws.Range("N135").Value = -70920.0

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H15").Value = 3000.0
$ws.Range("I15").Value = 0.0
$ws.Range("K15").Value = 0.0
$ws.Range("M15").ClearContents()
$ws.Range("H39").Value = 4949.5
$ws.Range("I39").Value = 4949.5
$ws.Range("K39").Value = 4949.5
$ws.Range("M39").Value = -4558.5
$ws.Range("H49").Value = 4949.5
$ws.Range("I49").Value = 4949.5
$ws.Range("K49").Value = 4949.5
$ws.Range("M49").Value = -4767.5
$ws.Range("I99").Value = 3187.5
$ws.Range("K99").Value = 3187.5
$ws.Range("M99").Value = -1689.5
$ws.Range("H111").Value = 75000.0
$ws.Range("J111").Value = 75000.0
$ws.Range("L111").Value = 75000.0
$ws.Range("N111").Value = -83180.0
$ws.Range("H118").Value = 0.0
$ws.Range("J118").Value = 0.0
$ws.Range("L118").Value = 0.0
$ws.Range("N118").ClearContents()
$ws.Range("I126").Value = 3187.5
$ws.Range("K126").Value = 9562.5
$ws.Range("M126").Value = -7092.5
$ws.Range("H132").Value = 3250.3635
$ws.Range("I132").Value = 3290.0527
$ws.Range("K132").Value = 9870.1581
$ws.Range("M132").Value = -7340.158100000001
$ws.Range("H134").Value = 4647.6665
$ws.Range("I134").Value = 4777.2
$ws.Range("J134").Value = 4000.0
$ws.Range("K134").Value = 14331.6
$ws.Range("L134").Value = 12000.0
$ws.Range("M134").Value = -11796.6
$ws.Range("N134").Value = -17070.0

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H16").Value = 714.0
$ws.Range("J16").Value = 2985.0
$ws.Range("L16").Value = 8955.0
$ws.Range("N16").Value = -9301.0
$ws.Range("H33").Value = 363.18182
$ws.Range("J33").Value = 490.0
$ws.Range("L33").Value = 2940.0
$ws.Range("N33").Value = -3506.0
$ws.Range("H80").Value = 0.0
$ws.Range("J80").Value = 0.0
$ws.Range("L80").Value = 0.0
$ws.Range("N80").ClearContents()
$ws.Range("H81").Value = 5437.5
$ws.Range("I81").Value = 750.0
$ws.Range("J81").Value = 7000.0
$ws.Range("K81").Value = 2250.0
$ws.Range("L81").Value = 21000.0
$ws.Range("M81").Value = -1127.0
$ws.Range("N81").Value = -23246.0
$ws.Range("H83").Value = 0.0
$ws.Range("J83").Value = 0.0
$ws.Range("L83").Value = 0.0
$ws.Range("N83").ClearContents()
$ws.Range("H84").Value = 5437.5
$ws.Range("I84").Value = 750.0
$ws.Range("J84").Value = 7000.0
$ws.Range("K84").Value = 6750.0
$ws.Range("L84").Value = 63000.0
$ws.Range("M84").Value = -1134.0
$ws.Range("N84").Value = -74232.0
$ws.Range("H88").Value = 12571.429
$ws.Range("J88").Value = 17111.111
$ws.Range("L88").Value = 51333.333
$ws.Range("N88").Value = -52189.333
$ws.Range("H91").Value = 12571.429
$ws.Range("J91").Value = 17111.111
$ws.Range("L91").Value = 51333.333
$ws.Range("N91").Value = -54297.333
$ws.Range("H117").Value = 1005.5
$ws.Range("I117").Value = 1037.3334
$ws.Range("J117").Value = 910.0
$ws.Range("K117").Value = 3112.0002
$ws.Range("L117").Value = 2730.0
$ws.Range("M117").Value = 329.9998000000001
$ws.Range("N117").Value = -9614.0
$ws.Range("H121").Value = 1419.6666
$ws.Range("I121").Value = 1543.3334
$ws.Range("K121").Value = 4630.0002
$ws.Range("M121").Value = -3320.0002

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 7444.923
$ws.Range("I122").Value = 7938.45
$ws.Range("K122").Value = 23815.35
$ws.Range("M122").Value = -21365.35
$ws.Range("H123").Value = 23570.428
$ws.Range("J123").Value = 23570.428
$ws.Range("L123").Value = 23570.428
$ws.Range("N123").Value = -28470.428
$ws.Range("H132").Value = 4354.1665
$ws.Range("I132").Value = 4354.1665
$ws.Range("K132").Value = 13062.4995
$ws.Range("M132").Value = -10532.4995

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5785.0
$ws.Range("I7").Value = 5740.0
$ws.Range("K7").Value = 5740.0
$ws.Range("M7").Value = -5628.0
$ws.Range("H14").Value = 13249.75
$ws.Range("J14").Value = 19500.0
$ws.Range("L14").Value = 19500.0
$ws.Range("N14").Value = -19844.0
$ws.Range("H16").Value = 1741.3846
$ws.Range("I16").Value = 919.5714
$ws.Range("J16").Value = 2700.1667
$ws.Range("K16").Value = 919.5714
$ws.Range("L16").Value = 2700.1667
$ws.Range("M16").Value = -749.5714
$ws.Range("N16").Value = -3040.1667
$ws.Range("H22").Value = 2922.4546
$ws.Range("J22").Value = 3205.9443
$ws.Range("L22").Value = 3205.9443
$ws.Range("N22").Value = -3795.9443
$ws.Range("H27").Value = 2922.4546
$ws.Range("J27").Value = 3205.9443
$ws.Range("L27").Value = 3205.9443
$ws.Range("N27").Value = -3419.9443
$ws.Range("H39").Value = 0.0
$ws.Range("I39").Value = 0.0
$ws.Range("K39").Value = 0.0
$ws.Range("M39").ClearContents()
$ws.Range("H40").Value = 7930.3184
$ws.Range("I40").Value = 5133.364
$ws.Range("K40").Value = 5133.364
$ws.Range("M40").Value = -4997.364
$ws.Range("H46").Value = 1737.1666
$ws.Range("I46").Value = 766.6667
$ws.Range("J46").Value = 1931.2667
$ws.Range("K46").Value = 766.6667
$ws.Range("L46").Value = 1931.2667
$ws.Range("M46").Value = -578.6667
$ws.Range("N46").Value = -2307.2667
$ws.Range("H47").Value = 0.0
$ws.Range("J47").Value = 0.0
$ws.Range("L47").Value = 0.0
$ws.Range("N47").ClearContents()
$ws.Range("H52").Value = 0.0
$ws.Range("J52").Value = 0.0
$ws.Range("L52").Value = 0.0
$ws.Range("N52").ClearContents()
$ws.Range("H53").Value = 23799.8
$ws.Range("I53").Value = 17249.75
$ws.Range("J53").Value = 50000.0
$ws.Range("K53").Value = 17249.75
$ws.Range("L53").Value = 50000.0
$ws.Range("M53").Value = -16731.75
$ws.Range("N53").Value = -51036.0
$ws.Range("H54").Value = 0.0
$ws.Range("J54").Value = 0.0
$ws.Range("L54").Value = 0.0
$ws.Range("N54").ClearContents()
$ws.Range("H55").Value = 1288.5416
$ws.Range("I55").Value = 415.4
$ws.Range("K55").Value = 415.4
$ws.Range("M55").Value = -242.4
$ws.Range("H82").Value = 4800.125
$ws.Range("I82").Value = 3380.2
$ws.Range("K82").Value = 3380.2
$ws.Range("M82").Value = -3019.2
$ws.Range("H85").Value = 4800.125
$ws.Range("I85").Value = 3380.2
$ws.Range("K85").Value = 3380.2
$ws.Range("M85").Value = -2132.2
$ws.Range("H126").Value = 5785.0
$ws.Range("I126").Value = 5740.0
$ws.Range("K126").Value = 17220.0
$ws.Range("M126").Value = -14750.0
$ws.Range("H132").Value = 4485.1704
$ws.Range("I132").Value = 4413.5454
$ws.Range("J132").Value = 4548.2
$ws.Range("K132").Value = 13240.6362
$ws.Range("L132").Value = 13644.6
$ws.Range("M132").Value = -10710.6362
$ws.Range("N132").Value = -18704.6

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H42").Value = 30000.0
$ws.Range("I42").Value = 30000.0
$ws.Range("J42").Value = 0.0
$ws.Range("K42").Value = 30000.0
$ws.Range("L42").Value = 0.0
$ws.Range("M42").Value = -29622.0
$ws.Range("N42").ClearContents()
$ws.Range("H46").Value = 94408.0
$ws.Range("J46").Value = 94408.0
$ws.Range("L46").Value = 94408.0
$ws.Range("N46").Value = -94870.0
$ws.Range("H62").Value = 8713.741
$ws.Range("J62").Value = 9389.637
$ws.Range("L62").Value = 9389.637
$ws.Range("N62").Value = -10637.637
$ws.Range("H65").Value = 8713.741
$ws.Range("J65").Value = 9389.637
$ws.Range("L65").Value = 46948.185
$ws.Range("N65").Value = -53188.185
$ws.Range("H97").Value = 25000.0
$ws.Range("J97").Value = 25000.0
$ws.Range("L97").Value = 25000.0
$ws.Range("N97").Value = -26982.0
$ws.Range("H107").Value = 2529.889
$ws.Range("I107").Value = 2009.8572
$ws.Range("K107").Value = 6029.571599999999
$ws.Range("M107").Value = -4109.571599999999
$ws.Range("H124").Value = 67294.664
$ws.Range("J124").Value = 73247.0
$ws.Range("L124").Value = 73247.0
$ws.Range("N124").Value = -83067.0
$ws.Range("H126").Value = 4624.647
$ws.Range("I126").Value = 4979.5454
$ws.Range("K126").Value = 14938.6362
$ws.Range("M126").Value = -12468.6362
$ws.Range("H132").Value = 5029.524
$ws.Range("I132").Value = 4531.0
$ws.Range("K132").Value = 13593.0
$ws.Range("M132").Value = -11063.0
$ws.Range("H134").Value = 94408.0
$ws.Range("J134").Value = 94408.0
$ws.Range("L134").Value = 283224.0
$ws.Range("N134").Value = -288294.0
